$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "41.694.20"
$ws.Range("E2").Value = "  +0.34%  "
Set-TextCell "D3" "2.229.87"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextCell "D5" "231.23"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextCell "D6" "0.622"
$ws.Range("E6").Value = "  -2.39%  "
Set-TextCell "D7" "60.09"
$ws.Range("E7").Value = "  -7.40%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextCell "D9" "0.403"
$ws.Range("E9").Value = "  -1.94%  "
Set-TextCell "D10" "57.92"
$ws.Range("E10").Value = "  -2.30%  "
Set-TextCell "D11" "0.0897"
$ws.Range("E11").Value = "  -0.20%  "
Set-TextCell "D12" "0.103"
$ws.Range("E12").Value = "  -1.06%  "
Set-TextCell "D13" "2.562.48"
$ws.Range("E13").Value = "  -0.93%  "
Set-TextCell "D14" "15.44"
$ws.Range("E14").Value = "  -4.87%  "
Set-TextCell "D15" "22.45"
$ws.Range("E15").Value = "  -0.66%  "
Set-TextCell "D16" "5.64"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  -4.39%  "
Set-TextCell "D18" "2.245.71"
$ws.Range("E18").Value = "  -0.54%  "
Set-TextCell "D19" "41.681.11"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -1.26%  "
Set-TextCell "D21" "72.42"
$ws.Range("E21").Value = "  -2.22%  "
Set-TextCell "D22" "6.13"
$ws.Range("E22").Value = "  -1.31%  "
Set-TextCell "D23" "247.54"
$ws.Range("E23").Value = "  -2.31%  "
Set-TextCell "D24" "0.999"
$ws.Range("E25").Value = "  -1.96%  "
Set-TextCell "D26" "2.30"
$ws.Range("E26").Value = "  -0.82%  "
Set-TextCell "D27" "9.78"
$ws.Range("E27").Value = "  -0.46%  "
Set-TextCell "D28" "169.33"
$ws.Range("E28").Value = "  -2.24%  "
Set-TextCell "D29" "0.141"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("E31").Value = "  -2.12%  "
Set-TextCell "D32" "2.57"
$ws.Range("E32").Value = "  -8.52%  "
$ws.Range("E33").Value = "  -2.02%  "
Set-TextCell "D34" "4.99"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("E35").Value = "  -1.57%  "
Set-TextCell "D36" "0.0654"
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  -9.10%  "
$ws.Range("E38").Value = "  -2.37%  "
Set-TextCell "D39" "3.59"
$ws.Range("E39").Value = "  -8.70%  "
Set-TextCell "D40" "0.000241"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +0.53%  "
Set-TextCell "D43" "8.62"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("E44").Value = "  -1.24%  "
Set-TextCell "D45" "98.85"
$ws.Range("E45").Value = "  -3.44%  "
Set-TextCell "D46" "0.0961"
$ws.Range("E46").Value = "  +1.82%  "
Set-TextCell "D47" "1.476.69"
$ws.Range("E47").Value = "  -2.66%  "
Set-TextCell "D48" "4.37"
$ws.Range("E48").Value = "  -10.33%  "
Set-TextCell "D49" "16.53"
$ws.Range("E49").Value = "  -8.14%  "
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -3.17%  "
